# Daily attendance processing - normalize "Recorded By" (column G) value ordering.
# Swaps the ordering of the recorder names/emails for a known set of
# placeholder/system values, leaving every other cell untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$colG = 7

$updated = 0

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colG)
    $val = $cell.Value2

    if ($val -eq "System, backup@backdoor.com, system") {
        $cell.Value2 = "system, System, backup@backdoor.com"
        $updated = $updated + 1
    }
    elseif ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
        $updated = $updated + 1
    }
    elseif ($val -eq "backup@backdoor.com, System") {
        $cell.Value2 = "System, backup@backdoor.com"
        $updated = $updated + 1
    }
}

Write-Host "Updated $updated cells in column G"
